# 1948_2019_Inflation_historic_Australia.xlsx
# - turn the hard-coded quarterly-average column (F) into real
#   =AVERAGE(Bn:En) formulas for every data row (2-72)
# - extend the summary block at the bottom (median) with new
#   average / max / min rows, and extend all of them to F72
# - rename the "ann" header to "ann  (ave quarters)"
# - cosmetic: column F width, scroll position / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column F: replace literal values with AVERAGE formulas -----------
# Row 2 through 72 (row 73/74 are blank spacer rows in the sheet).
$ws.Range("F2:F72").Formula = "=AVERAGE(B2:E2)"

# Make sure every cell in that range keeps/gets the 0.00% number format
# (a handful of rows previously used a plain 0% format).
$ws.Range("F2:F72").NumberFormat = "0.00%"

# --- summary block below the data --------------------------------------
$ws.Range("E75").Value = "median"
$ws.Range("F75").Formula = "=MEDIAN(F2:F72)"

$ws.Range("E76").Value = "average"
$ws.Range("F76").Formula = "=AVERAGE(F2:F72)"

$ws.Range("E77").Value = "max"
$ws.Range("F77").Formula = "=MAX(F2:F72)"

$ws.Range("E78").Value = "min"
$ws.Range("F78").Formula = "=MIN(F2:F72)"

$ws.Range("F75:F78").NumberFormat = "0.00%"

# --- header -------------------------------------------------------------
$ws.Range("F1").Value = "ann  (ave quarters)"

# --- cosmetics ----------------------------------------------------------
$ws.Columns("F").ColumnWidth = 17.85546875

$null = $ws.Range("F79").Select()
